$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 379, pushing the existing rows
# 379-433 down to 381-435 (data + formatting carried along automatically).
$ws.Rows("379:380").Insert()

# New row 379: weekly "Primera" quality entry for 2022-08-03.
$ws.Range("A379").Value = 8
$ws.Range("B379").Value = "Terminal La Palmera de La Serena"
$ws.Range("C379").Value = "Coquimbo"
$ws.Range("D379").Value = 44776
$ws.Range("E379").Value = 4
$ws.Range("F379").Value = "Fruta"
$ws.Range("G379").Value = 100101
$ws.Range("H379").Value = "Berries"
$ws.Range("I379").Value = 100101007
$ws.Range("J379").Value = "Kiwi"
$ws.Range("K379").Value = "Hayward"
$ws.Range("L379").Value = "Primera"
$ws.Range("M379").Value = 10
$ws.Range("N379").Value = 210000
$ws.Range("O379").Value = 220000
$ws.Range("P379").Value = 215000
$ws.Range("Q379").Value = "`$/bins (450 kilos)"
$ws.Range("R379").Value = "Región de O'Higgins"
$ws.Range("S379").Value = 478
$ws.Range("T379").Value = 450

# New row 380: weekly "Segunda" quality entry for 2022-08-03.
$ws.Range("A380").Value = 8
$ws.Range("B380").Value = "Terminal La Palmera de La Serena"
$ws.Range("C380").Value = "Coquimbo"
$ws.Range("D380").Value = 44776
$ws.Range("E380").Value = 4
$ws.Range("F380").Value = "Fruta"
$ws.Range("G380").Value = 100101
$ws.Range("H380").Value = "Berries"
$ws.Range("I380").Value = 100101007
$ws.Range("J380").Value = "Kiwi"
$ws.Range("K380").Value = "Hayward"
$ws.Range("L380").Value = "Segunda"
$ws.Range("M380").Value = 16
$ws.Range("N380").Value = 180000
$ws.Range("O380").Value = 190000
$ws.Range("P380").Value = 185000
$ws.Range("Q380").Value = "`$/bins (450 kilos)"
$ws.Range("R380").Value = "Región de O'Higgins"
$ws.Range("S380").Value = 411
$ws.Range("T380").Value = 450
